$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns before column D (shifts D:K -> F:M for the whole sheet,
#    including row heights/dimension/used-range; Excel keeps values & styles attached
#    to their rows as it pushes the old columns to the right).
$ws.Range("D1:E1").EntireColumn.Insert()

# 2) The newly inserted D:E columns inherit column C's plain format; restore the
#    correct number formats (date format for header rows, #,##0 for data rows) by
#    copying the format now sitting in column F/G (the old D/E) back onto D/E.
$ws.Range("F7:F102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Range("G7:G102").Copy()
$ws.Range("E7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new D:E columns with the newest two quarters of data (values taken
#    from the updated financial feed). Each entry is (row, D-value, E-value).
$newQuarterData = @(
    @(7, 43465, 43373),
    @(8, 530700, 533000),
    @(9, 201700, 198300),
    @(10, 329000, 334700),
    @(12, 7200, 6500),
    @(13, 0, 0),
    @(14, 65200, "NA"),
    @(15, 153600, 150600),
    @(17, 550200, 462900),
    @(18, -19500, 70100),
    @(20, -24500, 21300),
    @(21, 109600, 242000),
    @(22, 62200, 62100),
    @(23, -106200, 29300),
    @(24, 6400, 12800),
    @(25, 0, 0),
    @(26, -112600, 16500),
    @(27, -113200, 16100),
    @(28, 0, 0),
    @(29, 1000, 0),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, 24500, -21300),
    @(33, -112200, 16100),
    @(34, 0, 0),
    @(35, -112200, 16100),
    @(38, 43465, 43373),
    @(41, 928300, 1597300),
    @(42, 2282200, 1842900),
    @(43, 215300, 228600),
    @(44, 75400, 70100),
    @(45, 79700, 73400),
    @(46, 3580800, 3812200),
    @(47, 357600, 255200),
    @(48, 3415200, 3442500),
    @(49, 548100, 552100),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 759500, 788500),
    @(53, 0, 0),
    @(54, 8661300, 8850400),
    @(57, 123100, 120600),
    @(58, 959600, 1030100),
    @(59, 262400, 262200),
    @(60, 1345100, 1412900),
    @(61, 2573200, 2582800),
    @(62, 587500, 596000),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 4521100, 4606400),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 694100, 808300),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 4140200, 4244000),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, -112200, 16100),
    @(83, 153600, 150600),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 190600, 190400),
    @(91, -139900, -167200),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -744900, -376700),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, -115500, -7900),
    @(101, 1200, -1500),
    @(102, -668600, -195800)
)

foreach ($entry in $newQuarterData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 4).Value = $entry[1]
    $ws.Cells.Item($r, 5).Value = $entry[2]
}

# 4) Row 91 ("Capital Expenditures") was restated across the whole history, not just
#    shifted - overwrite D91:M91 outright with the corrected series.
$capexRow = @(-139900, -167200, -119600, -128500, -160600, -192100, -128100, -226600, -193900, -162500)
for ($i = 0; $i -lt $capexRow.Length; $i++) {
    $ws.Cells.Item(91, 4 + $i).Value = $capexRow[$i]
}

# 5) Two small restatements ripple into the historical total rows: "Net Borrowings"
#    (row 94) and "Change In Cash and Cash Equivalents" (row 102) both get their
#    2016-12-31 quarter (now column I) corrected by $100.
$ws.Cells.Item(94, 9).Value = -389000
$ws.Cells.Item(102, 9).Value = -186700
